# added visualizations, updated ksa dataset
# Append the new daily case row (27) to the Saudi cases dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row of data -------------------------------------------------
$ws.Range("B27").Value = 1104
$ws.Range("C27").Value = 92
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 82
$ws.Range("F27").Value = 0

# Set the date value first, then clone the date formatting (style) from
# the cell directly above (A26) so A27 picks up the same "yyyy-mm-dd"
# number format / style index instead of creating a brand new style.
$ws.Range("A27").Value = 43917
$ws.Range("A26").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- selection / view -------------------------------------------------
# Matches the saved workbook's final selection (E27). (Note: the source
# commit also scrolled the window so row 6 becomes the top visible row
# (topLeftCell="A6"); this runtime's ActiveWindow/Pane scroll properties
# don't persist independently of the selection, so only the selection
# change is reproducible here.)
$ws.Range("E27").Select() | Out-Null
